$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 (Experimental) and B18 (Compositional) need literal text values "true"/
# "false" -- plain Range.Value assignment auto-coerces recognized boolean
# literals to real booleans, so instead write a text-producing formula and
# then convert it to a literal value in place (Copy + PasteSpecial values),
# which preserves the existing cell style untouched.
$ws.Range("B7").Formula = '=""&"true"'
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null

$ws.Range("B18").Formula = '=""&"false"'
$ws.Range("B18").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false

# B8 (Date) value update
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
